$wb = $excel.ActiveWorkbook

# Both "NBR" and "BAR" sheets had their first 4 data rows (the rows where
# the Cutoff index A was 0-3) dropped. Excel's row-delete semantics shift
# everything below up by 4, which is exactly what the diff shows for
# columns B (Cutoff) and C (Reaction_number): the values that used to sit
# at B=5..19 now occupy rows 2-16. Column A is then renumbered back to a
# plain 0-based row index (0..14).
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A2:A5").EntireRow.Delete()

    for ($i = 0; $i -lt 15; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
